# for #DDDX 2012 conf
$wb = $excel.ActiveWorkbook

# --- Orders_schema: rename/retype the OrderId column, add a new PlacedOn/datetime row ---
$wsOrders = $wb.Worksheets.Item("Orders_schema")
$wsOrders.Range("A2").Value = "OrderNum"
$wsOrders.Range("B2").Value = "int"
$wsOrders.Range("A8").Value = "PlacedOn"
$wsOrders.Range("B8").Value = "datetime"

# --- OrderAdditions_schema: same OrderId -> OrderNum/int rename ---
$wsOrderAdditions = $wb.Worksheets.Item("OrderAdditions_schema")
$wsOrderAdditions.Range("A2").Value = "OrderNum"
$wsOrderAdditions.Range("B2").Value = "int"

# --- Products_data: fix the "Cappacinno" typo (preserve B3's quote-prefix style) ---
$wsProductsData = $wb.Worksheets.Item("Products_data")
$wsProductsData.Range("B3").Value = "Cappuccino"
$wsProductsData.Range("A3").Copy() | Out-Null
$wsProductsData.Range("B3").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# --- Selections on each sheet, then activate OrderAdditions_schema last so it becomes the tab shown ---
$wsProductsSchema = $wb.Worksheets.Item("Products_schema")
$wsProductsSchema.Activate()
$wsProductsSchema.Range("B2").Select() | Out-Null

$wsOrders.Activate()
$wsOrders.Range("B2").Select() | Out-Null

$wsProductsData.Activate()
$wsProductsData.Range("B4").Select() | Out-Null

$wsOrderAdditions.Activate()
$wsOrderAdditions.Range("B2").Select() | Out-Null
